# Applies the "Implicit, Explicit wait java files uploaded" edit to testData.xlsx:
#  - Sheet1: turns the password columns into plain numbers, drops the
#    "Mobile1" column, removes the demo data row (old row 4), and blanks
#    out (but keeps the styling of) the trailing date rows.
#  - Sheet2: keeps "amar"/"nath" (their shared-string ids shift down once
#    the now-unused strings are dropped).
#  - Trims the selection on Sheet1 to A1:D3.
#  - Shrinks the saved window width.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 : row 1 (user1/pass1/Mobile1/TRUE -> user1/12/TRUE) ---
$ws1.Cells.Item(1, 2).Value = 12
$ws1.Cells.Item(1, 3).ClearContents()

# --- Sheet1 : row 2 (user2/pass2/Mobile1/FALSE -> user2/23/FALSE) ---
$ws1.Cells.Item(2, 2).Value = 23
$ws1.Cells.Item(2, 3).ClearContents()

# --- Sheet1 : row 3 (user3/pass3/Mobile1/TRUE -> user3/34/TRUE) ---
$ws1.Cells.Item(3, 2).Value = 34
$ws1.Cells.Item(3, 3).ClearContents()

# --- Sheet1 : drop the old row 4 (123456/45678/amar/FALSE) entirely,
# then reinsert a blank row so the rows below keep their original
# row numbers (5 and 6). ---
$ws1.Rows.Item(4).Delete() | Out-Null
$ws1.Rows.Item(4).Insert() | Out-Null

# --- Sheet1 : blank the values of the old date rows (5 and 6) while
# keeping their number-format styles. ---
$ws1.Cells.Item(5, 1).ClearContents()
$ws1.Cells.Item(5, 2).ClearContents()
$ws1.Cells.Item(5, 3).ClearContents()
$ws1.Cells.Item(5, 4).ClearContents()
$ws1.Cells.Item(6, 1).ClearContents()

# --- Sheet1 : update the selected range shown when the sheet is opened ---
$ws1.Range("A1:D3").Select() | Out-Null

# --- Sheet2 : A1/B1 stay "amar"/"nath" (their shared-string index moves
# once the unused pass/Mobile1/date strings are gone). ---
$ws2.Range("A1").Value = "amar"
$ws2.Range("B1").Value = "nath"

# --- Workbook : shrink the saved window width ---
$wb.Windows.Item(1).Width = 13845

$wb.Save()
